# Applies updated percentage values to the "Female" (rows 3-6) and
# "Male" (rows 19-21) age-group blocks of the demographics table.
# Values are stored as plain text in the sheet (inline strings), so a
# leading apostrophe is used to force text entry and avoid Excel's
# automatic numeric coercion (which would otherwise introduce floating
# point artifacts such as 0.48999999999999999 instead of 0.49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Female / Age 60 - 69 (row 3) ---
$ws.Range("C3").Value = "'0.49"
$ws.Range("E3").Value = "'0.49"
$ws.Range("G3").Value = "'0.49"
$ws.Range("I3").Value = "'0.49"

# --- Female / Age 70 - 79 (row 4) ---
$ws.Range("C4").Value = "'0.29"
$ws.Range("E4").Value = "'0.29"
$ws.Range("G4").Value = "'0.29"
$ws.Range("I4").Value = "'0.29"

# --- Female / Age 80 - 89 (row 5) ---
$ws.Range("C5").Value = "'0.18"
$ws.Range("E5").Value = "'0.19"
$ws.Range("G5").Value = "'0.17"
$ws.Range("I5").Value = "'0.18"

# --- Female / Age 90 plus (row 6) ---
$ws.Range("C6").Value = "'0.04"
$ws.Range("E6").Value = "'0.04"
$ws.Range("G6").Value = "'0.06"
$ws.Range("I6").Value = "'0.05"

# --- Male / Age 60 - 69 (row 19) ---
$ws.Range("E19").Value = "'0.55"
$ws.Range("G19").Value = "'0.55"
$ws.Range("I19").Value = "'0.55"

# --- Male / Age 70 - 79 (row 20) ---
$ws.Range("C20").Value = "'0.29"
$ws.Range("E20").Value = "'0.29"
$ws.Range("G20").Value = "'0.29"
$ws.Range("I20").Value = "'0.29"

# --- Male / Age 80 - 89 (row 21) ---
$ws.Range("C21").Value = "'0.14"
$ws.Range("E21").Value = "'0.14"
$ws.Range("G21").Value = "'0.13"
